# Apply the BOM edit: add a new row 12 (mirroring the blank row 11 layout),
# extend the G-column multiplication formula down to G12,
# and update the workbook's active selection accordingly.
#
# Row 11 is blank (A:F empty, G11 = F11*E11 via the shared formula group
# that starts at G3). We replicate that same row layout/styling onto row 12
# by copying the whole row, then give G12 its own formula (copying a shared
# formula cell only carries over its cached value, not the formula itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:G11").Copy($ws.Range("A12:G12"))
$ws.Range("G12").Formula = "=F12*E12"

# Update selection to match the authored change (active cell moves to F12).
$ws.Range("F12").Select()
